$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Typography sheet: add a new typography entry "t_AI" on row 8
# ---------------------------------------------------------------
$tWs = $wb.Worksheets.Item("Typography")

# Make sure the new cells pick up no explicit cell-level style
# (same "Normal" / default style as the other data rows in this table).
$tWs.Range("B8:J8").Style = "Normal"

$tWs.Range("B8").Value = "t_AI"
$tWs.Range("C8").Value = "Asap-Regular.ttf"
$tWs.Range("D8").Value = 15
$tWs.Range("E8").Value = 4
$tWs.Range("F8").Value = "?"
$tWs.Range("H8").Value = "0-9"

# ---------------------------------------------------------------
# Translation sheet: point several rows at the new "t_AI" typography
# and update a couple of single-use text ids / values
# ---------------------------------------------------------------
$xWs = $wb.Worksheets.Item("Translation")

$xWs.Range("C15").Value = "t_AI"
$xWs.Range("C16").Value = "t_AI"
$xWs.Range("C17").Value = "t_AI"
$xWs.Range("C18").Value = "t_AI"
$xWs.Range("C19").Value = "t_AI"
$xWs.Range("C20").Value = "t_AI"
$xWs.Range("C21").Value = "t_AI"
$xWs.Range("C22").Value = "t_AI"

$xWs.Range("B27").Value = "SingleUseId83"
$xWs.Range("E27").Value = "<value> "

$xWs.Range("B28").Value = "SingleUseId84"
$xWs.Range("E28").Value = "0"
